$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: a new song row, reusing the same icon image + download host as
# row 2 but with its own song name.
$ws.Range("A3").Value = "https://cdn-icons-png.flaticon.com/512/2815/2815428.png"
$ws.Range("B3").Value = "Henlo I am longggggggggggggg"
$ws.Range("C3").Value = "https://www.pagalworld.com.sb/files/download/type/64/id/70390"

# Turn A3/C3 into real hyperlinks (adds the relationship + <hyperlink> entry).
$ws.Hyperlinks.Add($ws.Range("A3"), "https://cdn-icons-png.flaticon.com/512/2815/2815428.png")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.pagalworld.com.sb/files/download/type/64/id/70390")

# Hyperlinks.Add re-stamps the cell with the built-in "Hyperlink" style;
# re-apply the formatting copied from row 2 so A3/C3 end up styled exactly
# like A2/C2 (A2 uses a manual blue-underline font, C2 uses the named
# "Hyperlink" cell style) instead of a freshly minted duplicate style.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

$ws.Range("D3").Select() | Out-Null
